$wb = $excel.ActiveWorkbook

# --- "John" sheet: shift every year in column A down by one year ---
# (A2:A37 currently hold 2020..2055  ->  2021..2056)
$john = $wb.Worksheets.Item("John")
for ($r = 2; $r -le 37; $r++) {
    $cell = $john.Cells.Item($r, 1)
    $cell.Value = $cell.Value() + 1
}

# --- "Sally" sheet: shift every year in column A down by one year ---
# (A2:A40 currently hold 2020..2058  ->  2021..2059)
$sally = $wb.Worksheets.Item("Sally")
for ($r = 2; $r -le 40; $r++) {
    $cell = $sally.Cells.Item($r, 1)
    $cell.Value = $cell.Value() + 1
}

# --- Update the selection/active-cell on each sheet and which tab is active ---
# "John" ends up with A38 selected (its view had scrolled further down the list).
$john.Activate()
$null = $john.Range("A38").Select()

# "Sally" becomes the active (visible) tab, with B8 selected.
$sally.Activate()
$null = $sally.Range("B8").Select()
